# Applies the "Updated cryptos list on Thu Nov 14 17:55:12 UTC 2024 with GitHub
# Actions" commit: refreshes the Price (D) / Volume(1h) (E) figures for almost
# every coin row, and swaps row 51 from ImmutableX to ARBITRUM (name, link,
# price, volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A number of "Price" cells hold plain decimal-looking text (e.g. "215.52").
# The source data stores these as text, but typing a bare numeric-looking
# string into a normally-formatted cell makes Excel coerce it into a real
# number -- silently corrupting it (dropped trailing zeros, float rounding,
# scientific notation, ...). Temporarily mark those cells as Text so the
# assignment below is stored verbatim, then clear the temporary formatting
# again so the cell keeps its original (default) style, same as before.
$textPriceCells = @("D5", "D6", "D7", "D8", "D14", "D17", "D19", "D21", "D22", "D25", "D26", "D27", "D31", "D32", "D33", "D34", "D35", "D37", "D40", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Price (D) / Volume(1h) (E) refresh, row by row ---

$ws.Range("D2").Value = "89.255.73"
$ws.Range("E2").Value = "  -3.07%  "

$ws.Range("D3").Value = "3.137.65"
$ws.Range("E3").Value = "  -4.12%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "215.52"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").Value = "635.19"
$ws.Range("E6").Value = "  +0.93%  "

$ws.Range("D7").Value = "0.396"
$ws.Range("E7").Value = "  -3.69%  "

$ws.Range("D8").Value = "0.767"
$ws.Range("E8").Value = "  +6.96%  "

$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").Value = "3.135.61"
$ws.Range("E10").Value = "  -4.05%  "

$ws.Range("E11").Value = "  -5.26%  "

$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("E13").Value = "  -5.75%  "

$ws.Range("D14").Value = "5.32"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("D15").Value = "89.060.70"
$ws.Range("E15").Value = "  -3.06%  "

$ws.Range("D16").Value = "3.713.56"
$ws.Range("E16").Value = "  -3.96%  "

$ws.Range("D17").Value = "32.43"
$ws.Range("E17").Value = "  -5.53%  "

$ws.Range("D18").Value = "3.151.86"
$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("D19").Value = "0.0000232"
$ws.Range("E19").Value = "  +18.59%  "

$ws.Range("E20").Value = "  +1.13%  "

$ws.Range("D21").Value = "13.32"
$ws.Range("E21").Value = "  -5.05%  "

$ws.Range("D22").Value = "427.88"
$ws.Range("E22").Value = "  -2.66%  "

$ws.Range("E23").Value = "  -5.72%  "

$ws.Range("E24").Value = "  -6.69%  "

$ws.Range("D25").Value = "5.47"
$ws.Range("E25").Value = "  +1.76%  "

$ws.Range("D26").Value = "83.04"
$ws.Range("E26").Value = "  +7.62%  "

$ws.Range("D27").Value = "11.58"
$ws.Range("E27").Value = "  -6.37%  "

$ws.Range("D28").Value = "3.302.97"
$ws.Range("E28").Value = "  -4.31%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  -11.54%  "

$ws.Range("D31").Value = "0.962"
$ws.Range("E31").Value = "  -4.08%  "

$ws.Range("D32").Value = "4.04"
$ws.Range("E32").Value = "  +10.27%  "

$ws.Range("D33").Value = "8.21"
$ws.Range("E33").Value = "  -6.54%  "

$ws.Range("D34").Value = "507.80"
$ws.Range("E34").Value = "  -8.57%  "

$ws.Range("D35").Value = "0.147"
$ws.Range("E35").Value = "  +12.82%  "

$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").Value = "1.31"
$ws.Range("E37").Value = "  +1.49%  "

$ws.Range("E38").Value = "  -4.49%  "

$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("D40").Value = "22.25"

$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("E43").Value = "  -6.84%  "

$ws.Range("E44").Value = "  -7.62%  "

$ws.Range("D45").Value = "145.61"
$ws.Range("E45").Value = "  -2.99%  "

$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("D47").Value = "43.76"
$ws.Range("E47").Value = "  -3.21%  "

$ws.Range("D48").Value = "164.93"
$ws.Range("E48").Value = "  -8.14%  "

$ws.Range("D49").Value = "0.721"
$ws.Range("E49").Value = "  -1.19%  "

$ws.Range("D50").Value = "24.41"
$ws.Range("E50").Value = "  -2.76%  "

$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "0.598"
$ws.Range("E51").Value = "  -5.85%  "

# Restore the original (default) formatting on the cells we text-forced above.
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).ClearFormats()
}
